$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34 (shifts existing rows 34-45 down to 35-46,
# inheriting the D/E/N/O date-format styles from the old row 34 automatically).
$ws.Rows("34:34").Insert()

# Fill in the new row 34 with the "Luna Ciccardi" reservation.
$ws.Range("A34").Value = "Luna Ciccardi"
$ws.Range("B34").Value = "Booking"

# Phone number must be stored as literal text (keeps the leading "+").
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "+393485814122"

$ws.Range("D34").Value = 45891
$ws.Range("E34").Value = 45893
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 261.62
$ws.Range("H34").Value = 215.04
$ws.Range("I34").Value = 46.58
$ws.Range("J34").Value = 17.8
$ws.Range("K34").Value = 2025
$ws.Range("L34").Value = 8

# M34 stays blank; N34/O34 stay blank with no date formatting applied (unlike
# the other data rows, this new row has no N/O style).
$ws.Range("N34:O34").ClearFormats()

# The phone number for Franziska Lindermeier (now row 44 after the shift) loses
# its leading "+" and becomes a plain number in the resaved workbook.
$ws.Range("C44").NumberFormat = "General"
$ws.Range("C44").Value = 4917661016719
